# support_peplib3.xlsx update (Adrian v1.8)
# Removes the blank "x" separator row that used to follow each of the
# 13 peptide groups on Sheet1 (rows 64, 72, 80, 88, 96, 104), shifting
# all subsequent data up and shrinking the used range from A1:B104 to
# A1:B98.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete from the bottom up so earlier row numbers stay valid.
$blankRows = @(104, 96, 88, 80, 72, 64)
foreach ($r in $blankRows) {
    $ws.Rows.Item($r).Delete() | Out-Null
}

# Leave the selection where the author last left it.
$ws.Range("G103").Select() | Out-Null
